$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 11007
$ws1.Range("F4").Value = 268
$ws1.Range("F5").Value = 1224
$ws1.Range("F6").Value = 1099
$ws1.Range("F7").Value = 850
$ws1.Range("F8").Value = 288
$ws1.Range("F10").Value = 1178
$ws1.Range("F12").Value = 153
$ws1.Range("F13").Value = 899
$ws1.Range("F15").Value = 2049
$ws1.Range("G15").Value = 60
$ws1.Range("F16").Value = 22
$ws1.Range("F17").Value = 1017
$ws1.Range("F18").Value = 840
$ws1.Range("F20").Value = 814
$ws1.Range("F21").Value = 927
$ws1.Range("F24").Value = 641
$ws1.Range("F25").Value = 665
$ws1.Range("F26").Value = 130
$ws1.Range("F27").Value = 359
$ws1.Range("F28").Value = 1025
$ws1.Range("F29").Value = 49
$ws1.Range("F30").Value = 502
$ws1.Range("F31").Value = 179
$ws1.Range("F32").Value = 256
$ws1.Range("F34").Value = 588
$ws1.Range("F35").Value = 1899
$ws1.Range("F36").Value = 396
$ws1.Range("F37").Value = 38
$ws1.Range("F38").Value = 1452
$ws1.Range("F39").Value = 411
$ws1.Range("F45").Value = 79
$ws1.Range("F48").Value = 11
$ws1.Range("F49").Value = 84

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 200
$ws2.Range("F12").Value = 172
$ws2.Range("F14").Value = 141
$ws2.Range("F17").Value = 2

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2184
$ws3.Range("F3").Value = 640
$ws3.Range("F4").Value = 580

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2184
$ws4.Range("F3").Value = 640
$ws4.Range("F5").Value = 11007
$ws4.Range("F6").Value = 268
$ws4.Range("F8").Value = 580
$ws4.Range("F9").Value = 1099
$ws4.Range("F10").Value = 200
$ws4.Range("F11").Value = 1178
$ws4.Range("F13").Value = 153
$ws4.Range("F14").Value = 899
$ws4.Range("F15").Value = 2049
$ws4.Range("G15").Value = 60
$ws4.Range("F16").Value = 22
$ws4.Range("F17").Value = 1017
$ws4.Range("F18").Value = 840
$ws4.Range("F20").Value = 814
$ws4.Range("F21").Value = 927
$ws4.Range("F25").Value = 641
$ws4.Range("F28").Value = 665
$ws4.Range("F29").Value = 130
$ws4.Range("F30").Value = 359
$ws4.Range("F31").Value = 1025
$ws4.Range("F33").Value = 49
$ws4.Range("F34").Value = 502
$ws4.Range("F35").Value = 179
$ws4.Range("F36").Value = 256
$ws4.Range("F38").Value = 38
$ws4.Range("F39").Value = 1452
$ws4.Range("F40").Value = 411
$ws4.Range("F46").Value = 79
$ws4.Range("F48").Value = 84
